# power play title fix
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the title in C4: "Power Play Monster Manual" -> "Power Play Magic Manual"
$ws.Range("C4").Value = "Power Play Magic Manual"

# Re-set C5 so the shared-string table order matches (Power Play Progress stays the same text)
$ws.Range("C5").Value = "Power Play Progress"

# Update the active selection to C4
$ws.Range("C4").Select()
